# Applies the Chocobo_Profits profit-calc refresh (scheduled runner update).
# For each changed cell we either set the new numeric value, or -- when the
# diff shows the <c> element disappearing/appearing entirely (blank profit
# cell because HQ/NQ split produced no value that run) -- clear/create it.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1195.6428
$ws.Range("I6").Value = 658.0909
$ws.Range("J6").Value = 3166.6667
$ws.Range("K6").Value = 1974.2727
$ws.Range("L6").Value = 9500.000100000001
$ws.Range("M6").Value = -1862.2727
$ws.Range("N6").Value = -9724.000100000001
$ws.Range("H18").Value = 458.4
$ws.Range("J18").Value = 525.5
$ws.Range("L18").Value = 525.5
$ws.Range("N18").Value = -1093.5
$ws.Range("H38").Value = 4190.231
$ws.Range("I38").Value = 154.6
$ws.Range("J38").Value = 6712.5
$ws.Range("K38").Value = 463.8
$ws.Range("L38").Value = 20137.5
$ws.Range("M38").Value = -91.79999999999995
$ws.Range("N38").Value = -20881.5
$ws.Range("H39").Value = 467.72726
$ws.Range("I39").Value = 135
$ws.Range("J39").Value = 1050
$ws.Range("K39").Value = 405
$ws.Range("L39").Value = 3150
$ws.Range("M39").Value = -109
$ws.Range("N39").Value = -3742
$ws.Range("H43").Value = 2007.1818
$ws.Range("I43").Value = 1190.5
$ws.Range("J43").Value = 2188.6667
$ws.Range("K43").Value = 1190.5
$ws.Range("L43").Value = 2188.6667
$ws.Range("M43").Value = -1121.5
$ws.Range("N43").Value = -2326.6667
$ws.Range("H87").Value = 22600.867
$ws.Range("J87").Value = 22600.867
$ws.Range("L87").Value = 22600.867
$ws.Range("N87").Value = -25096.867
$ws.Range("H90").Value = 22600.867
$ws.Range("J90").Value = 22600.867
$ws.Range("L90").Value = 67802.601
$ws.Range("N90").Value = -80282.601
$ws.Range("H98").Value = 2333.3142
$ws.Range("I98").Value = 815.8261
$ws.Range("J98").Value = 5241.8335
$ws.Range("K98").Value = 815.8261
$ws.Range("L98").Value = 5241.8335
$ws.Range("M98").Value = 682.1739
$ws.Range("N98").Value = -8237.833500000001
$ws.Range("H116").Value = 328343.3
$ws.Range("I116").Value = 590240.75
$ws.Range("K116").Value = 590240.75
$ws.Range("M116").Value = -586798.75
$ws.Range("H122").Value = 2333.3142
$ws.Range("I122").Value = 815.8261
$ws.Range("J122").Value = 5241.8335
$ws.Range("K122").Value = 2447.4783
$ws.Range("L122").Value = 15725.5005
$ws.Range("M122").Value = 2.521700000000237
$ws.Range("N122").Value = -20625.5005
$ws.Range("H129").Value = 902.8043
$ws.Range("J129").Value = 963.225
$ws.Range("L129").Value = 2889.675
$ws.Range("N129").Value = -12889.675
$ws.Range("H138").Value = 5531.15
$ws.Range("I138").Value = 895.6896400000001
$ws.Range("J138").Value = 7424.507
$ws.Range("K138").Value = 2687.06892
$ws.Range("L138").Value = 22273.521
$ws.Range("M138").Value = 2452.93108
$ws.Range("N138").Value = -32553.521
$ws.Range("H141").Value = 7666.4375
$ws.Range("I141").Value = 7867.533
$ws.Range("K141").Value = 23602.599
$ws.Range("M141").Value = -18422.599

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1871.9231
$ws.Range("I61").Value = 1203.2858
$ws.Range("J61").Value = 2652
$ws.Range("K61").Value = 1203.2858
$ws.Range("L61").Value = 2652
$ws.Range("M61").Value = -991.2858000000001
$ws.Range("N61").Value = -3076
$ws.Range("H74").Value = 3515.1177
$ws.Range("I74").Value = 4073.8262
$ws.Range("J74").Value = 2346.9092
$ws.Range("K74").Value = 4073.8262
$ws.Range("L74").Value = 2346.9092
$ws.Range("M74").Value = -3199.8262
$ws.Range("N74").Value = -4094.9092
$ws.Range("H77").Value = 3515.1177
$ws.Range("I77").Value = 4073.8262
$ws.Range("J77").Value = 2346.9092
$ws.Range("K77").Value = 20369.131
$ws.Range("L77").Value = 11734.546
$ws.Range("M77").Value = -16001.131
$ws.Range("N77").Value = -20470.546
$ws.Range("H122").Value = 1792.8966
$ws.Range("I122").Value = 1090.2941
$ws.Range("J122").Value = 2788.25
$ws.Range("K122").Value = 3270.8823
$ws.Range("L122").Value = 8364.75
$ws.Range("M122").Value = -820.8823000000002
$ws.Range("N122").Value = -13264.75
$ws.Range("H132").Value = 2199.75
$ws.Range("I132").Value = 940.6667
$ws.Range("J132").Value = 4466.1
$ws.Range("K132").Value = 2822.0001
$ws.Range("L132").Value = 13398.3
$ws.Range("M132").Value = -292.0001000000002
$ws.Range("N132").Value = -18458.3
$ws.Range("H136").Value = 1871.9231
$ws.Range("I136").Value = 1203.2858
$ws.Range("J136").Value = 2652
$ws.Range("K136").Value = 3609.8574
$ws.Range("L136").Value = 7956
$ws.Range("M136").Value = -1059.8574
$ws.Range("N136").Value = -13056

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41671230
$ws.Range("I31").Value = 1569.6
$ws.Range("J31").Value = 71435270
$ws.Range("K31").Value = 1569.6
$ws.Range("L31").Value = 71435270
$ws.Range("M31").Value = -1274.6
$ws.Range("N31").Value = -71435860
$ws.Range("H34").Value = 41671230
$ws.Range("I34").Value = 1569.6
$ws.Range("J34").Value = 71435270
$ws.Range("K34").Value = 1569.6
$ws.Range("L34").Value = 71435270
$ws.Range("M34").Value = -1367.6
$ws.Range("N34").Value = -71435674
$ws.Range("H58").Value = 1606.6364
$ws.Range("I58").Value = 1492.6428
$ws.Range("J58").Value = 2049.9443
$ws.Range("K58").Value = 1492.6428
$ws.Range("L58").Value = 2049.9443
$ws.Range("M58").Value = -1289.6428
$ws.Range("N58").Value = -2455.9443
$ws.Range("H87").Value = 20212.5
$ws.Range("J87").Value = 20212.5
$ws.Range("L87").Value = 20212.5
$ws.Range("N87").Value = -22584.5
$ws.Range("H90").Value = 20212.5
$ws.Range("J90").Value = 20212.5
$ws.Range("L90").Value = 60637.5
$ws.Range("N90").Value = -72493.5
$ws.Range("H99").Value = 9528414
$ws.Range("I99").Value = 25002840
$ws.Range("J99").Value = 5691.5386
$ws.Range("K99").Value = 25002840
$ws.Range("L99").Value = 5691.5386
$ws.Range("M99").Value = -25001342
$ws.Range("N99").Value = -8687.5386
$ws.Range("H126").Value = 9528414
$ws.Range("I126").Value = 25002840
$ws.Range("J126").Value = 5691.5386
$ws.Range("K126").Value = 75008520
$ws.Range("L126").Value = 17074.6158
$ws.Range("M126").Value = -75006050
$ws.Range("N126").Value = -22014.6158
$ws.Range("H134").Value = 4682.3237
$ws.Range("I134").Value = 5240.864
$ws.Range("J134").Value = 3658.3333
$ws.Range("K134").Value = 15722.592
$ws.Range("L134").Value = 10974.9999
$ws.Range("M134").Value = -13187.592
$ws.Range("N134").Value = -16044.9999
$ws.Range("H136").Value = 1606.6364
$ws.Range("I136").Value = 1492.6428
$ws.Range("J136").Value = 2049.9443
$ws.Range("K136").Value = 4477.928400000001
$ws.Range("L136").Value = 6149.8329
$ws.Range("M136").Value = -1927.928400000001
$ws.Range("N136").Value = -11249.8329

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6944.4443
$ws.Range("J80").Value = 7500
$ws.Range("L80").Value = 22500
$ws.Range("N80").Value = -24372
$ws.Range("H83").Value = 6944.4443
$ws.Range("J83").Value = 7500
$ws.Range("L83").Value = 67500
$ws.Range("N83").Value = -76860
$ws.Range("H92").Value = 612.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 612.5
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("M92").Value = 1837.5
$ws.Range("N92").Value = -4333.5
$ws.Range("H98").Value = 300.85715
$ws.Range("I98").Value = 268.66666
$ws.Range("J98").Value = 325
$ws.Range("K98").Value = 805.9999799999999
$ws.Range("L98").Value = 975
$ws.Range("M98").Value = 692.0000200000001
$ws.Range("N98").Value = -3971
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H131").Value = 9616237
$ws.Range("I131").Value = 100000290
$ws.Range("J131").Value = 912.8298
$ws.Range("K131").Value = 300000870
$ws.Range("L131").Value = 2738.4894
$ws.Range("M131").Value = -299995830
$ws.Range("N131").Value = -12818.4894
$ws.Range("H132").Value = 1837.5454
$ws.Range("I132").Value = 635.2
$ws.Range("J132").Value = 2839.5
$ws.Range("K132").Value = 5716.8
$ws.Range("L132").Value = 25555.5
$ws.Range("M132").Value = -3186.8
$ws.Range("N132").Value = -30615.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 171500130
$ws.Range("I14").Value = 171500130
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 171500130
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -171499962
$ws.Range("H102").Value = 2021
$ws.Range("I102").Value = 1373.1
$ws.Range("K102").Value = 1373.1
$ws.Range("M102").Value = 248.9000000000001
$ws.Range("H107").Value = 11112088
$ws.Range("I107").Value = 483.75
$ws.Range("J107").Value = 18519824
$ws.Range("K107").Value = 483.75
$ws.Range("L107").Value = 18519824
$ws.Range("M107").Value = 1436.25
$ws.Range("N107").Value = -18523664
$ws.Range("H126").Value = 1853.34
$ws.Range("I126").Value = 1867.0104
$ws.Range("K126").Value = 5601.031199999999
$ws.Range("M126").Value = -3131.031199999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 85866
$ws.Range("I22").Value = 201398.2
$ws.Range("K22").Value = 201398.2
$ws.Range("M22").Value = -201103.2
$ws.Range("H27").Value = 85866
$ws.Range("I27").Value = 201398.2
$ws.Range("K27").Value = 201398.2
$ws.Range("M27").Value = -201291.2
$ws.Range("H136").Value = 2603.5
$ws.Range("I136").Value = 1369.826
$ws.Range("J136").Value = 5756.222
$ws.Range("K136").Value = 4109.478
$ws.Range("L136").Value = 17268.666
$ws.Range("M136").Value = -1559.478
$ws.Range("N136").Value = -22368.666

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 8893.5
$ws.Range("J12").Value = 8893.5
$ws.Range("L12").Value = 8893.5
$ws.Range("N12").Value = -9177.5
$ws.Range("H107").Value = 667.3333
$ws.Range("I107").Value = 402
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1206
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 714
$ws.Range("N107").Value = -6240
